$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = "Talk"
$ws.Range("D6").Value = "Walk"
$ws.Range("D8").Value = "Laugh"

$ws.Range("C6").Value = 2
$ws.Range("E6").Value = 187
$ws.Range("F6").Value = 217

$ws.Range("C7").Value = 3
$ws.Range("E7").Value = 218
$ws.Range("F7").Value = 318

$ws.Range("C8").Value = 4
$ws.Range("E8").Value = 319
$ws.Range("F8").Value = 383

$ws.Range("F8").Select()
